# Populate Sheet1 with a small username/password/是否执行 sample data set.
# Row 1 holds the headers, row 2 holds a sample record.
# The column order in which distinct string values are first introduced
# controls the order they land in xl/sharedStrings.xml, so the headers that
# repeat a string used in row 2 ("是否执行") are written after row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"

$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = 123456
$ws.Range("C2").Value = "Y"
$ws.Range("D2").Value = "Y"

$ws.Range("C1").Value = "是否执行"
$ws.Range("D1").Value = "是否执行"

# admin / 123456 are stored with a text number format so the numeric
# password keeps its literal digits instead of being treated as a number.
$ws.Range("A2:B2").NumberFormat = "@"

# Leave the active cell on B2, matching the saved selection.
$ws.Range("B2").Select() | Out-Null

# Match the sheet's printed page setup (A4/Letter-ish "9" = A4, portrait).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
